$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("What if user would like to update the Other Thread results in the UI")
# moves up to row 5. Capture its contents/style before clearing.
$movedText = $ws.Range("B6").Text

# Clear old row 6 cells (A6, B6) since they move to row 5.
$ws.Range("A6").Clear()
$ws.Range("B6").Clear()

# Recreate row 5: A5 gets the yellow-fill style (copy from A4/A8), B5 gets the moved text.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("B5").Value = $movedText

# C4 gets "Done" text, centered horizontally; C5 stays empty but shares the same
# centered alignment style so the merge looks consistent.
$ws.Range("C4").Value = "Done"
$ws.Range("C4:C5").HorizontalAlignment = -4108   # xlCenter

# Merge C4:C5
$ws.Range("C4:C5").Merge() | Out-Null

# Update the active selection to C6 (matches the post-edit sheetView selection)
$ws.Range("C6").Select() | Out-Null
